$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.055.35"
$ws.Range("E2").Value = "  -2.39%  "
$ws.Range("D3").Value = "3.007.46"
$ws.Range("E3").Value = "  -2.16%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.77%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -2.59%  "
$ws.Range("D9").Value = "3.004.46"
$ws.Range("E9").Value = "  -2.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.148"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.50%  "
$ws.Range("E11").Value = "  -2.86%  "
$ws.Range("E12").Value = "  -2.60%  "
$ws.Range("E13").Value = "  -4.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.33%  "
$ws.Range("E15").Value = "  +1.85%  "
$ws.Range("D16").Value = "3.502.86"
$ws.Range("E16").Value = "  -2.12%  "
$ws.Range("D17").Value = "62.048.27"
$ws.Range("E17").Value = "  -2.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.33%  "
$ws.Range("D19").Value = "3.005.83"
$ws.Range("E19").Value = "  -2.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "460.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.42%  "
$ws.Range("E22").Value = "  -4.23%  "
$ws.Range("E23").Value = "  -1.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.09%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  -2.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.81%  "
$ws.Range("E32").Value = "  -6.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.96%  "
$ws.Range("E34").Value = "  -3.98%  "
$ws.Range("E35").Value = "  -4.60%  "
$ws.Range("D36").Value = "0.0₃0785"
$ws.Range("E36").Value = "  -4.96%  "
$ws.Range("E37").Value = "  -5.47%  "
$ws.Range("E38").Value = "  -5.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -11.44%  "
$ws.Range("E42").Value = "  -0.41%  "
$ws.Range("E43").Value = "  -7.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0352"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "379.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -14.77%  "
$ws.Range("D46").Value = "2.746.38"
$ws.Range("E46").Value = "  -2.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "37.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.61%  "
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.70%  "
